# "Path to Graduation" planner — consolidate all terms' courses into the
# first ("2022") block and drop the now-empty "Fall 2024 / Spring 2024 /
# Summer 2024" block entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the third term block (rows 21-29: "Fall/Spring/Summer 2024"
# header + its Total row) completely — shifts nothing below it, so the
# sheet's used range shrinks from A2:F29 to A2:F20.
$ws.Rows("21:29").Delete()

# Consolidate the Spring/Summer 2022 column courses (previously spread
# across the 2022 and 2023 blocks) into rows 4-9 of the first block.
$ws.Range("C4").Value = "CPSC 4148"
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = "CPSC 4176"
$ws.Range("F4").Value = 3

$ws.Range("C5").Value = "CPSC 4155"
$ws.Range("D5").Value = 3

$ws.Range("C6").Value = "CPSC 4157"
$ws.Range("D6").Value = 3

$ws.Range("A7").Value = "CPSC 3165"
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = "CPSC 4175"
$ws.Range("D7").Value = 3

$ws.Range("A8").Value = "CPSC 4135"
$ws.Range("B8").Value = 3

$ws.Range("A9").Value = "CPSC 4000"
$ws.Range("B9").Value = 0

# The Fall 2023 block's detail rows are now redundant (their courses
# moved up into the first block) — clear them, leaving the "Fall 2023"
# header (row 12) and its Total row (row 20) untouched.
$ws.Range("A13:F15").ClearContents()
